# edit.ps1 - Applies the two changes described by the commit:
#   1. Update the cached "datetimeFigureOut" field text from 4/12/2023
#      to 4/24/2023 everywhere it appears (slide master + all slide layouts).
#   2. Fix the KDF reference on slide 3 from "NIST 800-56Ar2" to
#      "NIST 800-56Ar3".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder fields: 4/12/2023 -> 4/24/2023
# ---------------------------------------------------------------------
$oldDate = "4/12/2023"
$newDate = "4/24/2023"
$ppPlaceholderDate = 16

function Update-DatePlaceholders {
    param($shapes)

    $updated = 0
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)

        $isDatePlaceholder = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder -and $sh.HasTextFrame) {
            if ($sh.TextFrame.HasText) {
                if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                    $sh.TextFrame.TextRange.Text = $newDate
                    $updated++
                }
            }
        }
    }
    return $updated
}

$master = $p.SlideMaster

Update-DatePlaceholders($master.Shapes) | Out-Null

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-DatePlaceholders($layout.Shapes) | Out-Null
}

# ---------------------------------------------------------------------
# 2) KDF reference text fix on slide 3: NIST 800-56Ar2 -> NIST 800-56Ar3
# ---------------------------------------------------------------------
$oldRef = "NIST 800-56Ar2"
$newRef = "NIST 800-56Ar3"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
        $shp = $slide.Shapes.Item($k)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldRef) {
                    $shp.TextFrame.TextRange.Text = $newRef
                }
            }
        }
    }
}
